$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burndown Chart")

# Sprint backlog: log effort for day "Day 3" (column G) against
# "Preparation of the workplace" (row 6) and "Identifie code smells" (row 8).
$ws.Range("G6").Value = 1
$ws.Range("G8").Value = 0.4

# Formulas in row 11 (Completed Effort) and row 12 (Remaining Effort)
# recalculate automatically from the new inputs.

# Leave the active selection on N12, matching the saved view state.
$null = $ws.Range("N12").Select()
